$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptos price/volume snapshot.
# Cells that are pure numeric-looking text (e.g. "7.05", "1.00") are forced
# to Text format first so Excel keeps them as literal strings instead of
# silently converting them to numbers.
$ws.Range('D2').Value = '55.841.40'
$ws.Range('E2').Value = '  -3.99%  '
$ws.Range('D3').Value = '2.926.48'
$ws.Range('E3').Value = '  -4.47%  '
$ws.Range('E4').Value = '  +0.19%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '489.80'
$ws.Range('E5').Value = '  -7.05%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '131.96'
$ws.Range('E6').Value = '  -7.65%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.419'
$ws.Range('E8').Value = '  -6.56%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '7.05'
$ws.Range('E9').Value = '  -7.15%  '
$ws.Range('E10').Value = '  -8.29%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.346'
$ws.Range('E11').Value = '  -6.73%  '
$ws.Range('D12').Value = '3.435.51'
$ws.Range('E12').Value = '  -4.37%  '
$ws.Range('E13').Value = '  -4.03%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '25.63'
$ws.Range('E14').Value = '  -6.48%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000154'
$ws.Range('E15').Value = '  -10.97%  '
$ws.Range('D16').Value = '56.104.24'
$ws.Range('E16').Value = '  -3.48%  '
$ws.Range('D17').Value = '2.930.29'
$ws.Range('E17').Value = '  -4.60%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '5.90'
$ws.Range('E18').Value = '  -5.41%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.33'
$ws.Range('E19').Value = '  -6.34%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.63'
$ws.Range('E20').Value = '  -6.99%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '312.63'
$ws.Range('E21').Value = '  -8.71%  '
$ws.Range('E22').Value = '  -0.20%  '
$ws.Range('E23').Value = '  +0.24%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.476'
$ws.Range('E24').Value = '  -5.93%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '61.97'
$ws.Range('E25').Value = '  -5.22%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.159'
$ws.Range('E27').Value = '  -5.42%  '
$ws.Range('E28').Value = '  -13.91%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.34'
$ws.Range('E29').Value = '  -9.44%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.89'
$ws.Range('E30').Value = '  -8.33%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.72'
$ws.Range('E31').Value = '  -7.46%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '19.74'
$ws.Range('E32').Value = '  -6.83%  '
$ws.Range('E33').Value = '  -10.06%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '149.95'
$ws.Range('E34').Value = '  -5.32%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.39'
$ws.Range('E35').Value = '  -8.72%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.57'
$ws.Range('E36').Value = '  -6.77%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.18'
$ws.Range('E37').Value = '  -10.82%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '23.21'
$ws.Range('E38').Value = '  -11.34%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0643'
$ws.Range('E39').Value = '  -8.13%  '
$ws.Range('B40').Value = 'RenzoRestakedETH'
$ws.Range('C40').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D40').Value = '2.954.45'
$ws.Range('E40').Value = '  -4.69%  '
$ws.Range('B41').Value = 'OKB'
$ws.Range('C41').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '36.94'
$ws.Range('E41').Value = '  -2.36%  '
$ws.Range('E42').Value = '  +0.15%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.62'
$ws.Range('E43').Value = '  -8.57%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.629'
$ws.Range('E44').Value = '  -5.58%  '
$ws.Range('D45').Value = '2.116.07'
$ws.Range('E45').Value = '  -9.63%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.31'
$ws.Range('E46').Value = '  -11.23%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '5.78'
$ws.Range('E47').Value = '  -4.83%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.897'
$ws.Range('E48').Value = '  -13.47%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0226'
$ws.Range('E49').Value = '  -7.21%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '18.57'
$ws.Range('E50').Value = '  -8.13%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0837'
$ws.Range('E51').Value = '  -7.42%  '
